$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(344, 44418, 1, 8, 79.44389275074478),
    @(345, 44419, 0, 8, 79.44389275074478),
    @(346, 44420, 5, 10, 99.30486593843098),
    @(347, 44421, 0, 10, 99.30486593843098),
    @(348, 44422, 0, 10, 99.30486593843098),
    @(349, 44423, 0, 10, 99.30486593843098),
    @(350, 44424, 6, 12, 119.1658391261172),
    @(351, 44425, 0, 11, 109.2353525322741),
    @(352, 44426, 0, 11, 109.2353525322741),
    @(353, 44427, 1, 7, 69.51340615690168),
    @(354, 44428, 4, 11, 109.2353525322741),
    @(355, 44429, 1, 12, 119.1658391261172),
    @(356, 44430, 0, 12, 119.1658391261172),
    @(357, 44431, 1, 7, 69.51340615690168)
)

# Copy the formatting of the last existing data row (343) into column A
# of each new row so the date number format / style matches the rest
# of the column.
$ws.Cells.Item(343, 1).Copy() | Out-Null

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
